# Update of league bases, 30-05-2024 23:16
# 1) Rename header columns I1/J1 (ht_goals_h/ht_goals_a -> HTHG/HTAG)
# 2) Rotate the match-odds data among rows 176-178 (id/date columns A & D
#    stay put; everything else for each row takes on the values that used
#    to belong to the next row, wrapping 176 <- 177 <- 178 <- 176).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header rename ---------------------------------------------------
$ws.Range("I1").Value = "HTHG"
$ws.Range("J1").Value = "HTAG"

# --- 2. Rotate rows 176-178 ----------------------------------------------

# Row 176
$ws.Range("B176").Value = 7217242
$ws.Range("C176").Value = "Azerbaijan Premier League"
$ws.Range("E176").Value = "Zira IK"
$ws.Range("F176").Value = "Sabail FC"
$ws.Range("G176").Value = 1
$ws.Range("H176").Value = 0
$ws.Range("I176").Value = 0
$ws.Range("J176").Value = 0
$ws.Range("K176").Value = "H"
$ws.Range("L176").Value = 1.727
$ws.Range("M176").Value = 3.25
$ws.Range("N176").Value = 4.333
$ws.Range("O176").Value = 1.45
$ws.Range("P176").Value = 3.4
$ws.Range("Q176").Value = 7
$ws.Range("R176").Value = -1
$ws.Range("S176").Value = 1.8
$ws.Range("T176").Value = 2
$ws.Range("U176").Value = 2.5
$ws.Range("V176").Value = 1.95
$ws.Range("W176").Value = 1.75
$ws.Range("X176").Value = 0.45
$ws.Range("Y176").Value = -1
$ws.Range("Z176").Value = -1
$ws.Range("AA176").Value = 0
$ws.Range("AB176").Value = 0
$ws.Range("AC176").Value = -1
$ws.Range("AD176").Value = 0.75

# Row 177
$ws.Range("B177").Value = 7217874
$ws.Range("C177").Value = "Azerbaijan Premier League"
$ws.Range("E177").Value = "FK Sumqayit"
$ws.Range("F177").Value = "FK Gabala"
$ws.Range("G177").Value = 1
$ws.Range("H177").Value = 0
$ws.Range("I177").Value = 0
$ws.Range("J177").Value = 0
$ws.Range("K177").Value = "H"
$ws.Range("L177").Value = 1.333
$ws.Range("M177").Value = 4.5
$ws.Range("N177").Value = 7
$ws.Range("O177").Value = 1.55
$ws.Range("P177").Value = 3.9
$ws.Range("Q177").Value = 4.5
$ws.Range("R177").Value = -1
$ws.Range("S177").Value = 1.975
$ws.Range("T177").Value = 1.825
$ws.Range("U177").Value = 2.5
$ws.Range("V177").Value = 1.75
$ws.Range("W177").Value = 1.95
$ws.Range("X177").Value = 0.55
$ws.Range("Y177").Value = -1
$ws.Range("Z177").Value = -1
$ws.Range("AA177").Value = 0
$ws.Range("AB177").Value = 0
$ws.Range("AC177").Value = -1
$ws.Range("AD177").Value = 0.95

# Row 178
$ws.Range("B178").Value = 7223244
$ws.Range("C178").Value = "Azerbaijan Premier League"
$ws.Range("E178").Value = "Neftchi Baku"
$ws.Range("F178").Value = "Sabah"
$ws.Range("G178").Value = 0
$ws.Range("H178").Value = 1
$ws.Range("I178").Value = 0
$ws.Range("J178").Value = 0
$ws.Range("K178").Value = "A"
$ws.Range("L178").Value = 2.1
$ws.Range("M178").Value = 3
$ws.Range("N178").Value = 3.25
$ws.Range("O178").Value = 1.8
$ws.Range("P178").Value = 3.2
$ws.Range("Q178").Value = 4
$ws.Range("R178").Value = -0.5
$ws.Range("S178").Value = 1.825
$ws.Range("T178").Value = 1.975
$ws.Range("U178").Value = 2.75
$ws.Range("V178").Value = 1.95
$ws.Range("W178").Value = 1.85
$ws.Range("X178").Value = -1
$ws.Range("Y178").Value = -1
$ws.Range("Z178").Value = 3
$ws.Range("AA178").Value = -1
$ws.Range("AB178").Value = 0.9750000000000001
$ws.Range("AC178").Value = -1
$ws.Range("AD178").Value = 0.8500000000000001
